$d = $word.ActiveDocument

function Find-ParagraphContaining($doc, $needle) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) "Major design decisions" body: replace the placeholder sentence with the
#    real write-up about remodeling the ScreenHandler to use an observer
#    pattern.
# ---------------------------------------------------------------------------
$oldDesign = "No major design decisions have been made this week."
$newDesign = "We remodeled the ScreenHandler to implement a observer pattern. The reason for this is that  the observer pattern is a tried and tested method which feels more reliable then our original homebrewed code."
$d.Content.Find.Execute($oldDesign, $true, $false, $false, $false, $false, $true, 1, $false, $newDesign, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Add a body paragraph under "Issues, problems and risks" (it previously
#    had no text under the heading).
# ---------------------------------------------------------------------------
$issuesHeading = Find-ParagraphContaining $d "Issues, problems and risks"
$issuesHeading.Range.InsertParagraphAfter() | Out-Null
$issuesBody = $issuesHeading.Next()
$issuesBody.Range.Text = "The group has not found any current risks, issues or problems."
$issuesBody.Style = "Normal"

# ---------------------------------------------------------------------------
# 3) Add a body paragraph under "Current Status" (it previously had no text
#    under the heading either).
# ---------------------------------------------------------------------------
$statusHeading = Find-ParagraphContaining $d "Current Status"
$statusHeading.Range.InsertParagraphAfter() | Out-Null
$statusBody = $statusHeading.Next()
$statusBody.Range.Text = "The basic coding for the project is now complete and can hopefully start directly with implementing the game world and gameplay after the L2 delivery."
$statusBody.Style = "Normal"

# ---------------------------------------------------------------------------
# 4) "Planned work" body: expand the sentence about updating the
#    architecture document and add the bit about continuing to code the
#    basic project structure.
# ---------------------------------------------------------------------------
$oldPlanned = "Make updates to our architecture, Finish the UML class diagram and the L2 documentation. Might also try to do one or two more UML state machine diagrams."
$newPlanned = "Make updates to our architecture document that we handed in with the L1, Finish the UML class diagram and the L2 documentation. Might also try to do one or two more UML state machine diagrams as well as continue coding the basic structure for our project."
$d.Content.Find.Execute($oldPlanned, $true, $false, $false, $false, $false, $true, 1, $false, $newPlanned, 2) | Out-Null
